# Updated cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '57.780.88'
$ws.Range("E2").Value = '  -1.03%  '

# Row 3
$ws.Range("D3").Value = '2.446.39'

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.41'
$ws.Range("E5").Value = '  +0.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.89'
$ws.Range("E6").Value = '  -2.16%  '

# Row 7
$ws.Range("E7").Value = '  +0.47%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.561'
$ws.Range("E8").Value = '  +0.64%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0975'
$ws.Range("E9").Value = '  +0.06%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.92'
$ws.Range("E11").Value = '  -4.35%  '

# Row 12
$ws.Range("E12").Value = '  -3.19%  '

# Row 13
$ws.Range("D13").Value = '2.880.47'
$ws.Range("E13").Value = '  -2.59%  '

# Row 14
$ws.Range("D14").Value = '57.679.54'
$ws.Range("E14").Value = '  -1.19%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.70'
$ws.Range("E15").Value = '  -1.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000133'
$ws.Range("E16").Value = '  -1.89%  '

# Row 17
$ws.Range("D17").Value = '2.449.96'
$ws.Range("E17").Value = '  -2.19%  '

# Row 18
$ws.Range("E18").Value = '  -2.87%  '

# Row 19
$ws.Range("E19").Value = '  -0.81%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '314.53'
$ws.Range("E20").Value = '  -2.37%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.10'
$ws.Range("E21").Value = '  +0.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.00%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.87'
$ws.Range("E23").Value = '  +0.73%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.409'
$ws.Range("E24").Value = '  +1.95%  '

# Row 25
$ws.Range("E25").Value = '  +0.59%  '

# Row 26
$ws.Range("E26").Value = '  -2.64%  '

# Row 27
$ws.Range("E27").Value = '  -2.36%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.91'
$ws.Range("E28").Value = '  +2.43%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0734'
$ws.Range("E29").Value = '  -2.79%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.70'
$ws.Range("E30").Value = '  -1.38%  '

# Row 31
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.10'
$ws.Range("E31").Value = '  -2.61%  '

# Row 32
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  -4.62%  '

# Row 33
$ws.Range("E33").Value = '  +0.02%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  +0.16%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.80'
$ws.Range("E35").Value = '  -1.61%  '

# Row 36
$ws.Range("E36").Value = '  -6.48%  '

# Row 37
$ws.Range("E37").Value = '  -4.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.28'
$ws.Range("E38").Value = '  +0.54%  '

# Row 39
$ws.Range("E39").Value = '  -1.15%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.790'
$ws.Range("E40").Value = '  +1.35%  '

# Row 41
$ws.Range("E41").Value = '  -2.91%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '264.54'
$ws.Range("E42").Value = '  -4.93%  '

# Row 43
$ws.Range("E43").Value = '  -2.85%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.80'
$ws.Range("E44").Value = '  -3.71%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '123.86'
$ws.Range("E45").Value = '  +0.12%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0928'
$ws.Range("E46").Value = '  +1.02%  '

# Row 47
$ws.Range("E47").Value = '  -1.94%  '

# Row 48
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0211'
$ws.Range("E48").Value = '  -1.54%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.01'
$ws.Range("E49").Value = '  -4.65%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.30'
$ws.Range("E50").Value = '  -3.80%  '

# Row 51
$ws.Range("D51").Value = '1.698.94'
$ws.Range("E51").Value = '  -2.65%  '

